$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record for Espinaca (Vega Monumental Concepción) was
# reported for date 2022-04-08 (serial 44659). Insert it as a new row 21,
# pushing all the following rows (previously 21-66) down by one (22-67).
$ws.Rows("21:21").Insert()

$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44659
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 100112012
$ws.Range("G21").Value = "Espinaca"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 6500
$ws.Range("M21").Value = 6200
$ws.Range("N21").Value = "$/cuna 10 kilos"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 620
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = "Hortaliza"
